## Dec20 - order object
## Adds a new "Order" worksheet between "Contract" and "LoginPage",
## populated with test-case data (TC_ID / AccountName / ContractNumber / Comments),
## and nudges the leftover cell-selection on "Contract" to match the saved file.

$wb = $excel.ActiveWorkbook

$wsContract = $wb.Worksheets.Item("Contract")
$wsLogin    = $wb.Worksheets.Item("LoginPage")

# Leave a "last selection" on Contract like the author did, without disturbing
# which sheet is active when the new sheet gets inserted below.
$wsContract.Range("D15").Select() | Out-Null
$wsLogin.Activate() | Out-Null

# LoginPage is the active sheet, so Add() inserts the new sheet right before
# it -> final order becomes Contract, Order, LoginPage.
$order = $wb.Worksheets.Add()
$order.Name = "Order"

# ---- Header row (bold) ----
$order.Range("A1:D1").Font.Bold = $true
$order.Range("C1").NumberFormat = "@"

$order.Range("A1").Value = "TC_ID"
$order.Range("B1").Value = "AccountName"
$order.Range("C1").Value = "ContractNumber"
$order.Range("D1").Value = "Comments"

# ---- ContractNumber column keeps its leading zeros / non-numeric text ----
$order.Range("C2:C5").NumberFormat = "@"

# ---- Data rows ----
$order.Range("A2").Value = "TC_001"
$order.Range("B2").Value = "Kevin Testing"
$order.Range("C2").Value = "00000101"
$order.Range("D2").Value = "Valid data"

$order.Range("A3").Value = "TC_002"
$order.Range("B3").Value = "Kirthy"
$order.Range("C3").Value = "333"
$order.Range("D3").Value = "Invalid data"

$order.Range("A4").Value = "TC_003"
$order.Range("B4").Value = "Kevin Testing"
$order.Range("C4").Value = "00000103"

$order.Range("A5").Value = "TC_004"
$order.Range("B5").Value = "Cathrine"
$order.Range("C5").Value = "00000102"
$order.Range("D5").Value = "Valid data"

$order.Range("F13").Select() | Out-Null
